# Generate Report for Handoff
# Updates the localization-status report:
#  - Overview sheet: "Latest HO Xliff Generate Date" for the handoff batch
#    moves from 06:25:30 to 06:26:05
#  - zh-cn sheet: "Latest Handoff Datetime" moves from 06:25:22 to 06:25:55,
#    and the newly-handed-off rows get a "ht" Priority
#  - de-de sheet: "Latest Handoff Datetime" moves from 06:25:30 to 06:26:05,
#    and the newly-handed-off rows get a "ht" Priority

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 12, 13)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-13 06:26:05"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-13 06:25:55"
    $wsZhCn.Range("E$r").Value = "ht"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-13 06:26:05"
    $wsDeDe.Range("E$r").Value = "ht"
}
